$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q1" sheet right before the "总计" (totals) sheet.
#    Its layout mirrors the other quarterly sheets (2020-Q4 .. 2021-Q4):
#    header row in B1:H1, one data row in A2:H2.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q1")
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Pull over the exact header/column-A formatting (bold, bordered, centered)
# used by every other quarterly sheet so the new tab matches visually.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

# Fund code / name / percentages are stored as text (as in the other tabs),
# so force a text number format before assigning to avoid Excel coercing
# these numeric-looking strings into actual numbers.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "010551"
$newSheet.Range("C2").Value = "淳厚欣颐一年持有期混合"
$newSheet.Range("D2").Value = "3.41"
$newSheet.Range("E2").Value = "81.14"
$newSheet.Range("F2").Value = "2.71"
$newSheet.Range("G2").Value = "0.0924"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2. Add the matching "2022-Q1" summary row at the top of the "总计" sheet's
#    data, pushing the existing quarters down by one row.
#    (Re-fetch the sheet by name - inserting/renaming sheet tabs above can
#    shift positional handles, so grab a fresh reference here.)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-use the formatting of the (now shifted) first data row for the newly
# inserted one, so column A keeps its bold/bordered/centered style.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.09

# The running index in column A is a plain 0-based row counter - renumber it
# now that a row was inserted at the top.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
